$d = $word.ActiveDocument

# 1. Version number 1.0 -> 1.2.5
$d.Content.Find.Execute("1.0", $false, $false, $false, $false, $false, $true, 1, $false, "1.2.5", 2) | Out-Null

# 2. Creation -> Update
$d.Content.Find.Execute("Creation", $false, $false, $false, $false, $false, $true, 1, $false, "Update", 2) | Out-Null

# 3. Author name Fabrício Araújo -> Julio Paiva
$d.Content.Find.Execute("Fabrício Araújo", $false, $false, $false, $false, $false, $true, 1, $false, "Julio Paiva", 2) | Out-Null

# 4. Date 09/07/2020 -> 31/05/2023
$d.Content.Find.Execute("09/07/2020", $false, $false, $false, $false, $false, $true, 1, $false, "31/05/2023", 2) | Out-Null

# 5. Fix typo "usuario" -> "usuário" and add a period at end
$d.Content.Find.Execute("O usuario devidamente autenticado e na tela inicial do sistema", $false, $false, $false, $false, $false, $true, 1, $false, "O usuário devidamente autenticado e na tela inicial do sistema.", 2) | Out-Null

# 6. Fix typo "histório" -> "histórico"
$d.Content.Find.Execute("histório da tramitação", $false, $false, $false, $false, $false, $true, 1, $false, "histórico da tramitação", 2) | Out-Null

# 7. Add period after "Detalhar Diárias "
$d.Content.Find.Execute("2. System Apresenta a tela de Detalhar Diárias ", $false, $false, $false, $false, $false, $true, 1, $false, "2. System Apresenta a tela de Detalhar Diárias. ", 2) | Out-Null

# 8. Fix wording "Permite não permite" -> "Não permite"
$d.Content.Find.Execute("Permite não permite um novo envio", $false, $false, $false, $false, $false, $true, 1, $false, "Não permite um novo envio", 2) | Out-Null
